$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.351.84"
$ws.Range("E2").Value = '  -0.14%  '

$ws.Range("D3").Value = "'1.872.89"
$ws.Range("E3").Value = '  -1.06%  '

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = "'238.14"
$ws.Range("E5").Value = '  +0.34%  '

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("D7").Value = "'0.4784"
$ws.Range("E7").Value = '  -1.17%  '

$ws.Range("D8").Value = "'0.2822"
$ws.Range("E8").Value = '  -2.82%  '

$ws.Range("D9").Value = "'0.06510"
$ws.Range("E9").Value = '  -1.49%  '

$ws.Range("D10").Value = "'1.871.22"
$ws.Range("E10").Value = '  -2.22%  '

$ws.Range("D11").Value = "'0.07456"

$ws.Range("D12").Value = "'16.54"
$ws.Range("E12").Value = '  -2.34%  '

$ws.Range("D13").Value = "'5.106"
$ws.Range("E13").Value = '  -1.28%  '

$ws.Range("D14").Value = "'88.16"
$ws.Range("E14").Value = '  +0.41%  '

$ws.Range("D15").Value = "'0.6542"
$ws.Range("E15").Value = '  -1.20%  '

$ws.Range("D16").Value = "'30.333.15"
$ws.Range("E16").Value = '  -0.13%  '

$ws.Range("D17").Value = "'13.32"
$ws.Range("E17").Value = '  -1.03%  '

$ws.Range("D18").Value = "'0.9995"
$ws.Range("E18").Value = '  +0.17%  '

$ws.Range("D19").Value = "'0.000007587"
$ws.Range("E19").Value = '  -2.56%  '

$ws.Range("D20").Value = "'2.114.99"
$ws.Range("E20").Value = '  -1.11%  '

$ws.Range("D21").Value = "'5.299"
$ws.Range("E21").Value = '  -2.25%  '

$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").Value = "'219.20"
$ws.Range("E23").Value = '  +12.98%  '

$ws.Range("D24").Value = "'6.188"
$ws.Range("E24").Value = '  -0.11%  '

$ws.Range("D25").Value = "'9.330"
$ws.Range("E25").Value = '  -0.33%  '

$ws.Range("D26").Value = "'167.71"
$ws.Range("E26").Value = '  +1.74%  '

$ws.Range("D27").Value = "'18.44"
$ws.Range("E27").Value = '  +1.30%  '

$ws.Range("D28").Value = "'1.974"
$ws.Range("E28").Value = '  +1.76%  '

$ws.Range("D29").Value = "'1.447"
$ws.Range("E29").Value = '  -0.35%  '

$ws.Range("D30").Value = "'0.09356"
$ws.Range("E30").Value = '  +2.10%  '

$ws.Range("D31").Value = "'4.316"
$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("D32").Value = "'4.032"
$ws.Range("E32").Value = '  -0.46%  '

$ws.Range("D33").Value = "'0.05076"
$ws.Range("E33").Value = '  -0.31%  '

$ws.Range("D34").Value = "'1.206"
$ws.Range("E34").Value = '  +4.91%  '

$ws.Range("D35").Value = "'0.7527"
$ws.Range("E35").Value = '  +3.54%  '

$ws.Range("D36").Value = "'2.714"
$ws.Range("E36").Value = '  +0.60%  '

$ws.Range("D37").Value = "'0.01824"
$ws.Range("E37").Value = '  +1.74%  '

$ws.Range("D38").Value = "'2.612"
$ws.Range("E38").Value = '  -1.45%  '

$ws.Range("D39").Value = "'2.068"
$ws.Range("E39").Value = '  -0.47%  '

$ws.Range("D40").Value = "'0.9063"
$ws.Range("E40").Value = '  -1.70%  '

$ws.Range("D41").Value = "'106.83"
$ws.Range("E41").Value = '  +0.64%  '

$ws.Range("D42").Value = "'5.898"
$ws.Range("E42").Value = '  +0.49%  '

$ws.Range("D43").Value = "'0.4272"
$ws.Range("E43").Value = '  -1.11%  '

$ws.Range("D44").Value = "'1.004"
$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("D45").Value = "'7.379"
$ws.Range("E45").Value = '  -1.67%  '

$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = "'0.1284"
$ws.Range("E46").Value = '  -3.54%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = "'64.15"
$ws.Range("E47").Value = '  -1.47%  '

$ws.Range("D48").Value = "'8.934"
$ws.Range("E48").Value = '  -0.03%  '

$ws.Range("D49").Value = "'1.468"
$ws.Range("E49").Value = '  -7.06%  '

$ws.Range("D50").Value = "'33.59"
$ws.Range("E50").Value = '  -1.28%  '

$ws.Range("D51").Value = "'0.3889"
$ws.Range("E51").Value = '  +0.59%  '
